# Apply the "Generar estadisticas de tiempos" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update underlying data values (dependent formulas in J16, D17, J17, D18, E18, J18, K18
# will recalculate automatically).
$ws.Range("K16").Value = 933
$ws.Range("E17").Value = 82
$ws.Range("K17").Value = 1651

# Update sheet selection / view state (activeCell=I15, sqref=I15:L22).
$ws.Range("I15:L22").Select()
